$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# New elementalResistance (column O) entry for the Swarm row (row 6)
$ws.Range("O6").Value = "ICE:-0.2"

# New magicalResistance (column N) entries for several rows
$ws.Range("N11").Value = 0.15
$ws.Range("N13").Value = -0.2
$ws.Range("N14").Value = 0.5
$ws.Range("N32").Value = -0.25
$ws.Range("N36").Value = -0.15

# maxDamage for row 16 changed from 18 to 16
$ws.Range("J16").Value = 16

# elementalResistance text updates
$ws.Range("O20").Value = "SHADOW:0.2 POISON:-0.15"
$ws.Range("O28").Value = "SHADOW:0.2 POISON:0.15 FIRE:0.15 HOLY:-0.2 ICE:-0.15"

# Move the active selection from N18 to N9 (cursor position change recorded in the sheet view)
$ws.Range("N9").Select()
